$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 2712.8572  # H17: was 3498
$ws.Cells.Item(17, 10).Value = 2712.8572  # J17: was 3498
$ws.Cells.Item(17, 12).Value = 8138.571599999999  # L17: was 10494
$ws.Cells.Item(17, 14).Value = -8474.571599999999  # N17: was -10830
$ws.Cells.Item(92, 8).Value = 1002.8333  # H92: was 900.5714
$ws.Cells.Item(92, 10).Value = 1399.3334  # J92: was 954.4
$ws.Cells.Item(92, 12).Value = 1399.3334  # L92: was 954.4
$ws.Cells.Item(92, 14).Value = -3895.3334  # N92: was -3450.4
$ws.Cells.Item(98, 8).Value = 665  # H98: was 717.8823
$ws.Cells.Item(98, 9).Value = 760.8  # I98: was 820.6429000000001
$ws.Cells.Item(98, 10).Value = 377.6  # J98: was 238.33333
$ws.Cells.Item(98, 11).Value = 760.8  # K98: was 820.6429000000001
$ws.Cells.Item(98, 12).Value = 377.6  # L98: was 238.33333
$ws.Cells.Item(98, 13).Value = 737.2  # M98: was 677.3570999999999
$ws.Cells.Item(98, 14).Value = -3373.6  # N98: was -3234.33333
$ws.Cells.Item(107, 8).Value = 599.8  # H107: was 617.55554
$ws.Cells.Item(107, 9).Value = 642.8570999999999  # I107: was 676.6667
$ws.Cells.Item(107, 11).Value = 642.8570999999999  # K107: was 676.6667
$ws.Cells.Item(107, 13).Value = 1277.1429  # M107: was 1243.3333
$ws.Cells.Item(112, 8).Value = 1816.2  # H112: was 1799.1818
$ws.Cells.Item(112, 9).Value = 1375  # I112: was 1300
$ws.Cells.Item(112, 10).Value = 1926.5  # J112: was 1986.375
$ws.Cells.Item(112, 11).Value = 4125  # K112: was 3900
$ws.Cells.Item(112, 12).Value = 5779.5  # L112: was 5959.125
$ws.Cells.Item(112, 13).Value = -3017  # M112: was -2792
$ws.Cells.Item(112, 14).Value = -7995.5  # N112: was -8175.125
$ws.Cells.Item(113, 8).Value = 3495  # H113: was 3499
$ws.Cells.Item(113, 10).Value = 3495  # J113: was 3499
$ws.Cells.Item(113, 12).Value = 3495  # L113: was 3499
$ws.Cells.Item(113, 14).Value = -10003  # N113: was -10007
$ws.Cells.Item(122, 8).Value = 665  # H122: was 717.8823
$ws.Cells.Item(122, 9).Value = 760.8  # I122: was 820.6429000000001
$ws.Cells.Item(122, 10).Value = 377.6  # J122: was 238.33333
$ws.Cells.Item(122, 11).Value = 2282.4  # K122: was 2461.9287
$ws.Cells.Item(122, 12).Value = 1132.8  # L122: was 714.99999
$ws.Cells.Item(122, 13).Value = 167.6000000000004  # M122: was -11.92870000000039
$ws.Cells.Item(122, 14).Value = -6032.8  # N122: was -5614.99999
$ws.Cells.Item(125, 8).Value = 1881.125  # H125: was 1999.8334
$ws.Cells.Item(125, 9).Value = 1200  # I125: was 0
$ws.Cells.Item(125, 10).Value = 1978.4286  # J125: was 1999.8334
$ws.Cells.Item(125, 11).Value = 10800  # K125: was 0
$ws.Cells.Item(125, 12).Value = 17805.8574  # L125: was 17998.5006
$ws.Cells.Item(125, 13).Value = -8340  # M125: was None
$ws.Cells.Item(125, 14).Value = -22725.8574  # N125: was -22918.5006
$ws.Cells.Item(130, 8).Value = 88888  # H130: was 88887.5
$ws.Cells.Item(130, 10).Value = 88888  # J130: was 88887.5
$ws.Cells.Item(130, 12).Value = 88888  # L130: was 88887.5
$ws.Cells.Item(130, 14).Value = -98928  # N130: was -98927.5
$ws.Cells.Item(135, 8).Value = 1710  # H135: was 1495
$ws.Cells.Item(135, 10).Value = 3500  # J135: was 4000
$ws.Cells.Item(135, 12).Value = 31500  # L135: was 36000
$ws.Cells.Item(135, 14).Value = -36570  # N135: was -41070
$ws.Cells.Item(137, 8).Value = 2538.625  # H137: was 2500.7576
$ws.Cells.Item(137, 9).Value = 1499.9445  # I137: was 1479.35
$ws.Cells.Item(137, 10).Value = 3874.0715  # J137: was 4072.1538
$ws.Cells.Item(137, 11).Value = 4499.833500000001  # K137: was 4438.049999999999
$ws.Cells.Item(137, 12).Value = 11622.2145  # L137: was 12216.4614
$ws.Cells.Item(137, 13).Value = -1949.833500000001  # M137: was -1888.049999999999
$ws.Cells.Item(137, 14).Value = -16722.2145  # N137: was -17316.4614

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 7665.025  # H32: was 7835.875
$ws.Cells.Item(32, 9).Value = 7349.1167  # I32: was 7524.8027
$ws.Cells.Item(32, 11).Value = 7349.1167  # K32: was 7524.8027
$ws.Cells.Item(32, 13).Value = -7062.1167  # M32: was -7237.8027
$ws.Cells.Item(61, 8).Value = 1285.1842  # H61: was 1309.8379
$ws.Cells.Item(61, 9).Value = 932.86664  # I61: was 952.1724
$ws.Cells.Item(61, 11).Value = 932.86664  # K61: was 952.1724
$ws.Cells.Item(61, 13).Value = -720.86664  # M61: was -740.1724
$ws.Cells.Item(74, 8).Value = 1616.2667  # H74: was 1392.6578
$ws.Cells.Item(74, 9).Value = 1160.3478  # I74: was 1003.9032
$ws.Cells.Item(74, 11).Value = 1160.3478  # K74: was 1003.9032
$ws.Cells.Item(74, 13).Value = -286.3478  # M74: was -129.9032
$ws.Cells.Item(77, 8).Value = 1616.2667  # H77: was 1392.6578
$ws.Cells.Item(77, 9).Value = 1160.3478  # I77: was 1003.9032
$ws.Cells.Item(77, 11).Value = 5801.739  # K77: was 5019.516
$ws.Cells.Item(77, 13).Value = -1433.739  # M77: was -651.5159999999996
$ws.Cells.Item(110, 8).Value = 1123  # H110: was 1065.5
$ws.Cells.Item(110, 9).Value = 1081  # I110: was 1080.6
$ws.Cells.Item(110, 10).Value = 1333  # J110: was 990
$ws.Cells.Item(110, 11).Value = 1081  # K110: was 1080.6
$ws.Cells.Item(110, 12).Value = 1333  # L110: was 990
$ws.Cells.Item(110, 13).Value = 964  # M110: was 964.4000000000001
$ws.Cells.Item(110, 14).Value = -5423  # N110: was -5080
$ws.Cells.Item(132, 8).Value = 2778.5264  # H132: was 2902.4707
$ws.Cells.Item(132, 9).Value = 2233.25  # I132: was 2334.9
$ws.Cells.Item(132, 11).Value = 6699.75  # K132: was 7004.700000000001
$ws.Cells.Item(132, 13).Value = -4169.75  # M132: was -4474.700000000001
$ws.Cells.Item(136, 8).Value = 1285.1842  # H136: was 1309.8379
$ws.Cells.Item(136, 9).Value = 932.86664  # I136: was 952.1724
$ws.Cells.Item(136, 11).Value = 2798.59992  # K136: was 2856.5172
$ws.Cells.Item(136, 13).Value = -248.5999199999997  # M136: was -306.5172000000002

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 873.2222  # H94: was 954
$ws.Cells.Item(94, 9).Value = 565.75  # I94: was 614.1429000000001
$ws.Cells.Item(94, 11).Value = 565.75  # K94: was 614.1429000000001
$ws.Cells.Item(94, 13).Value = -114.75  # M94: was -163.1429000000001

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1367.8  # H16: was 1614
$ws.Cells.Item(16, 9).Value = 1279.6666  # I16: was 1842
$ws.Cells.Item(16, 11).Value = 1279.6666  # K16: was 1842
$ws.Cells.Item(16, 13).Value = -992.6666  # M16: was -1555
$ws.Cells.Item(31, 8).Value = 2463.6365  # H31: was 2872
$ws.Cells.Item(31, 9).Value = 2722.7144  # I31: was 2532.2
$ws.Cells.Item(31, 10).Value = 2010.25  # J31: was 4004.6667
$ws.Cells.Item(31, 11).Value = 2722.7144  # K31: was 2532.2
$ws.Cells.Item(31, 12).Value = 2010.25  # L31: was 4004.6667
$ws.Cells.Item(31, 13).Value = -2427.7144  # M31: was -2237.2
$ws.Cells.Item(31, 14).Value = -2600.25  # N31: was -4594.6667
$ws.Cells.Item(34, 8).Value = 2463.6365  # H34: was 2872
$ws.Cells.Item(34, 9).Value = 2722.7144  # I34: was 2532.2
$ws.Cells.Item(34, 10).Value = 2010.25  # J34: was 4004.6667
$ws.Cells.Item(34, 11).Value = 2722.7144  # K34: was 2532.2
$ws.Cells.Item(34, 12).Value = 2010.25  # L34: was 4004.6667
$ws.Cells.Item(34, 13).Value = -2520.7144  # M34: was -2330.2
$ws.Cells.Item(34, 14).Value = -2414.25  # N34: was -4408.6667
$ws.Cells.Item(58, 8).Value = 2691.3333  # H58: was 2912.4443
$ws.Cells.Item(58, 9).Value = 2626.8  # I58: was 2883.4285
$ws.Cells.Item(58, 11).Value = 2626.8  # K58: was 2883.4285
$ws.Cells.Item(58, 13).Value = -2423.8  # M58: was -2680.4285
$ws.Cells.Item(113, 8).Value = 1367.8  # H113: was 1614
$ws.Cells.Item(113, 9).Value = 1279.6666  # I113: was 1842
$ws.Cells.Item(113, 11).Value = 1279.6666  # K113: was 1842
$ws.Cells.Item(113, 13).Value = 890.3334  # M113: was 328
$ws.Cells.Item(132, 8).Value = 4495.091  # H132: was 4513.2
$ws.Cells.Item(132, 9).Value = 4383.1113  # I132: was 4391.75
$ws.Cells.Item(132, 11).Value = 13149.3339  # K132: was 13175.25
$ws.Cells.Item(132, 13).Value = -10619.3339  # M132: was -10645.25
$ws.Cells.Item(134, 8).Value = 2700.5557  # H134: was 2829.7144
$ws.Cells.Item(134, 9).Value = 2829.2856  # I134: was 2968
$ws.Cells.Item(134, 10).Value = 2250  # J134: was 2000
$ws.Cells.Item(134, 11).Value = 8487.856800000001  # K134: was 8904
$ws.Cells.Item(134, 12).Value = 6750  # L134: was 6000
$ws.Cells.Item(134, 13).Value = -5952.856800000001  # M134: was -6369
$ws.Cells.Item(134, 14).Value = -11820  # N134: was -11070
$ws.Cells.Item(136, 8).Value = 2691.3333  # H136: was 2912.4443
$ws.Cells.Item(136, 9).Value = 2626.8  # I136: was 2883.4285
$ws.Cells.Item(136, 11).Value = 7880.400000000001  # K136: was 8650.2855
$ws.Cells.Item(136, 13).Value = -5330.400000000001  # M136: was -6100.2855
$ws.Cells.Item(141, 8).Value = 51184.535  # H141: was 49740.266
$ws.Cells.Item(141, 9).Value = 0  # I141: was 20000
$ws.Cells.Item(141, 10).Value = 51184.535  # J141: was 51864.57
$ws.Cells.Item(141, 11).Value = 0  # K141: was 20000
$ws.Cells.Item(141, 12).Value = 51184.535  # L141: was 51864.57
$ws.Cells.Item(141, 13).ClearContents()  # M141: was -14820
$ws.Cells.Item(141, 14).Value = -61544.535  # N141: was -62224.57

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(23, 8).Value = 293.5  # H23: was 302.57144
$ws.Cells.Item(23, 10).Value = 293.5  # J23: was 302.57144
$ws.Cells.Item(23, 12).Value = 880.5  # L23: was 907.71432
$ws.Cells.Item(23, 14).Value = -1350.5  # N23: was -1377.71432
$ws.Cells.Item(131, 8).Value = 1732.6  # H131: was 1706.0714
$ws.Cells.Item(131, 9).Value = 1250.4  # I131: was 1258
$ws.Cells.Item(131, 10).Value = 1973.7  # J131: was 2042.125
$ws.Cells.Item(131, 11).Value = 3751.2  # K131: was 3774
$ws.Cells.Item(131, 12).Value = 5921.1  # L131: was 6126.375
$ws.Cells.Item(131, 13).Value = 1288.8  # M131: was 1266
$ws.Cells.Item(131, 14).Value = -16001.1  # N131: was -16206.375

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 65.7  # H2: was 69.28570999999999
$ws.Cells.Item(2, 9).Value = 61.88889  # I2: was 72.5
$ws.Cells.Item(2, 10).Value = 100  # J2: was 50
$ws.Cells.Item(2, 11).Value = 61.88889  # K2: was 72.5
$ws.Cells.Item(2, 12).Value = 100  # L2: was 50
$ws.Cells.Item(2, 13).Value = 51.11111  # M2: was 40.5
$ws.Cells.Item(2, 14).Value = -326  # N2: was -276
$ws.Cells.Item(45, 8).Value = 24999.666  # H45: was 40000
$ws.Cells.Item(45, 10).Value = 24999.666  # J45: was 40000
$ws.Cells.Item(45, 12).Value = 24999.666  # L45: was 40000
$ws.Cells.Item(45, 14).Value = -26117.666  # N45: was -41118
$ws.Cells.Item(97, 8).Value = 735.4545000000001  # H97: was 640.86957
$ws.Cells.Item(97, 9).Value = 680.58826  # I97: was 589.41174
$ws.Cells.Item(97, 10).Value = 922  # J97: was 786.6667
$ws.Cells.Item(97, 11).Value = 680.58826  # K97: was 589.41174
$ws.Cells.Item(97, 12).Value = 922  # L97: was 786.6667
$ws.Cells.Item(97, 13).Value = -184.58826  # M97: was -93.41174000000001
$ws.Cells.Item(97, 14).Value = -1914  # N97: was -1778.6667
$ws.Cells.Item(122, 8).Value = 8000  # H122: was 0
$ws.Cells.Item(122, 9).Value = 8000  # I122: was 0
$ws.Cells.Item(122, 11).Value = 24000  # K122: was 0
$ws.Cells.Item(122, 13).Value = -21550  # M122: was None
$ws.Cells.Item(126, 8).Value = 2898.5  # H126: was 2899
$ws.Cells.Item(126, 9).Value = 2898  # I126: was 0
$ws.Cells.Item(126, 11).Value = 8694  # K126: was 0
$ws.Cells.Item(126, 13).Value = -6224  # M126: was None
$ws.Cells.Item(132, 8).Value = 3802.1428  # H132: was 3869
$ws.Cells.Item(132, 9).Value = 3600.6191  # I132: was 3680.8
$ws.Cells.Item(132, 11).Value = 10801.8573  # K132: was 11042.4
$ws.Cells.Item(132, 13).Value = -8271.8573  # M132: was -8512.400000000001

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 10900  # H122: was 6040.6665
$ws.Cells.Item(122, 9).Value = 10900  # I122: was 6040.6665
$ws.Cells.Item(122, 11).Value = 32700  # K122: was 18121.9995
$ws.Cells.Item(122, 13).Value = -30250  # M122: was -15671.9995

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(96, 8).Value = 1399.3334  # H96: was 1296
$ws.Cells.Item(96, 9).Value = 1399.3334  # I96: was 1296
$ws.Cells.Item(96, 11).Value = 1399.3334  # K96: was 1296
$ws.Cells.Item(96, 13).Value = -26.33339999999998  # M96: was 77
$ws.Cells.Item(113, 8).Value = 7701  # H113: was 4968.3335
$ws.Cells.Item(113, 9).Value = 7701  # I113: was 13902
$ws.Cells.Item(113, 10).Value = 0  # J113: was 501.5
$ws.Cells.Item(113, 11).Value = 23103  # K113: was 41706
$ws.Cells.Item(113, 12).Value = 0  # L113: was 1504.5
$ws.Cells.Item(113, 13).Value = -20933  # M113: was -39536
$ws.Cells.Item(113, 14).ClearContents()  # N113: was -5844.5
$ws.Cells.Item(138, 8).Value = 125000  # H138: was 0
$ws.Cells.Item(138, 10).Value = 125000  # J138: was 0
$ws.Cells.Item(138, 12).Value = 125000  # L138: was 0
$ws.Cells.Item(138, 14).Value = -135280  # N138: was None
